# Refresh crypto price/volume snapshot (scheduled GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.431.75"
$ws.Range("E2").Value = "  -0.53%  "

# Row 3
$ws.Range("D3").Value = "1.803.25"
$ws.Range("E3").Value = "  +0.59%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'228.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "

# Row 6
$ws.Range("E6").Value = "  +4.97%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "'35.03"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.77%  "

# Row 9
$ws.Range("D9").Value = "'0.300"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.16%  "

# Row 10
$ws.Range("D10").Value = "'0.0693"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.14%  "

# Row 11
$ws.Range("E11").Value = "  +0.14%  "

# Row 12
$ws.Range("D12").Value = "2.062.45"
$ws.Range("E12").Value = "  +0.49%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'11.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.15%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.814.09"
$ws.Range("E14").Value = "  +0.79%  "

# Row 15
$ws.Range("D15").Value = "'0.643"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.17%  "

# Row 16
$ws.Range("D16").Value = "34.387.71"
$ws.Range("E16").Value = "  -0.49%  "

# Row 17
$ws.Range("D17").Value = "'4.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.66%  "

# Row 18
$ws.Range("D18").Value = "'69.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.47%  "

# Row 19
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'245.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.93%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0797"
$ws.Range("E20").Value = "  -0.25%  "

# Row 21
$ws.Range("D21").Value = "'11.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.48%  "

# Row 22
$ws.Range("E22").Value = "  +0.14%  "

# Row 23
$ws.Range("D23").Value = "'4.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.20%  "

# Row 24
$ws.Range("D24").Value = "'172.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.24%  "

# Row 25
$ws.Range("D25").Value = "'2.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.26%  "

# Row 26
$ws.Range("D26").Value = "'7.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.59%  "

# Row 27
$ws.Range("E27").Value = "  +2.99%  "

# Row 28
$ws.Range("D28").Value = "'16.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.29%  "

# Row 29
$ws.Range("E29").Value = "  -1.14%  "

# Row 30
$ws.Range("E30").Value = "  -3.63%  "

# Row 31
$ws.Range("D31").Value = "'0.0529"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.12%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.22%  "

# Row 33
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.72%  "

# Row 34
$ws.Range("E34").Value = "  +0.24%  "

# Row 35
$ws.Range("D35").Value = "1.396.26"
$ws.Range("E35").Value = "  -2.14%  "

# Row 36
$ws.Range("D36").Value = "'0.680"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.30%  "

# Row 37
$ws.Range("D37").Value = "'2.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.26%  "

# Row 38
$ws.Range("D38").Value = "'1.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.58%  "

# Row 39
$ws.Range("E39").Value = "  -0.88%  "

# Row 40
$ws.Range("D40").Value = "'83.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.22%  "

# Row 41
$ws.Range("E41").Value = "  +2.86%  "

# Row 42
$ws.Range("D42").Value = "'0.949"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.78%  "

# Row 43
$ws.Range("E43").Value = "  -0.63%  "

# Row 44
$ws.Range("D44").Value = "'13.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.20%  "

# Row 45
$ws.Range("E45").Value = "  +3.28%  "

# Row 46
$ws.Range("E46").Value = "  -3.63%  "

# Row 47
$ws.Range("E47").Value = "  -2.23%  "

# Row 48
$ws.Range("D48").Value = "1.963.08"
$ws.Range("E48").Value = "  +0.51%  "

# Row 49
$ws.Range("D49").Value = "'104.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.21%  "

# Row 50
$ws.Range("E50").Value = "  +0.04%  "

# Row 51
$ws.Range("E51").Value = "  +0.51%  "
